# Insert a new weekly data row at row 201 (pushing existing rows 201-273
# down to 202-274), then populate the new row with the latest observation.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(201).Insert()

$ws.Range("A201").Value = 3
$ws.Range("B201").Value = 'Femacal de La Calera'
$ws.Range("C201").Value = 'Coquimbo'
$ws.Range("D201").Value = 45009
$ws.Range("E201").Value = 5
$ws.Range("F201").Value = 100112030
$ws.Range("G201").Value = 'Poroto granado'
$ws.Range("H201").Value = 'Sin especificar'
$ws.Range("I201").Value = 'Primera'
$ws.Range("J201").Value = 65
$ws.Range("K201").Value = 30000
$ws.Range("L201").Value = 31000
$ws.Range("M201").Value = 30538
$ws.Range("N201").Value = '$/saco 25 kilos'
$ws.Range("O201").Value = 'Provincia de Quillota'
$ws.Range("P201").Value = 1222
$ws.Range("Q201").Value = 25
$ws.Range("R201").Value = 'Hortaliza'
